# Add two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# formatting of the existing header/data columns (e.g. H = "IP").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1): copy format from H1 (bold/border/center-top)
# onto I1:J1, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data cells (rows 2-38): plain numeric values, no special style
# (matches the rest of the C:H data columns).
$dataI = @(3,9,3,4,8,6,6,6,4,7,7,10,9,10,5,9,9,5,8,6,7,8,11,6,9,7,4,7,6,9,6,5,6,7,8,5,8)
$dataJ = @(4,9,3,5,8,6,6,7,4,8,8,10,9,10,6,9,9,6,9,6,7,9,11,6,9,7,5,7,7,9,6,5,6,7,8,5,8)

for ($k = 0; $k -lt $dataI.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$k]
    $ws.Cells.Item($row, 10).Value = $dataJ[$k]
}

Write-Output "I0/IF columns added"
